$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gphb5"
$ws.Range("C2").Value = "Tshr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7013873333333334
$ws.Range("H2").Value = 2.104162
$ws.Range("I2").Value = 0.1758510422341793
$ws.Range("J2").Value = 0.1758510422341793
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7209850000000001
$ws.Range("N2").Value = 2.162955
$ws.Range("O2").Value = 0.1533325535399077
$ws.Range("P2").Value = 0.1533325535399077
$ws.Range("Q2").Value = 0.5056897465233334
$ws.Range("R2").Value = 4.551207718710001
$ws.Range("S2").Value = 0.02696368934842087
$ws.Range("T2").Value = 0.02696368934842087

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gphb5"
$ws.Range("C3").Value = "Tshr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.7013873333333334
$ws.Range("H3").Value = 2.104162
$ws.Range("I3").Value = 0.1758510422341793
$ws.Range("J3").Value = 0.1758510422341793
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.387303666666666
$ws.Range("N3").Value = 10.161911
$ws.Range("O3").Value = 0.7203810354238886
$ws.Range("P3").Value = 0.7203810354238885
$ws.Range("Q3").Value = 2.375811885953556
$ws.Range("R3").Value = 21.382306973582
$ws.Range("S3").Value = 0.1266797558850281
$ws.Range("T3").Value = 0.126679755885028

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gphb5"
$ws.Range("C4").Value = "Tshr"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.7013873333333334
$ws.Range("H4").Value = 2.104162
$ws.Range("I4").Value = 0.1758510422341793
$ws.Range("J4").Value = 0.1758510422341793
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1167496666666667
$ws.Range("N4").Value = 0.350249
$ws.Range("O4").Value = 0.02482926068494219
$ws.Range("P4").Value = 0.02482926068494219
$ws.Range("Q4").Value = 0.0818867373708889
$ws.Range("R4").Value = 0.7369806363380002
$ws.Range("S4").Value = 0.004366251369351217
$ws.Range("T4").Value = 0.004366251369351216

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Gphb5"
$ws.Range("C5").Value = "Tshr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.7013873333333334
$ws.Range("H5").Value = 2.104162
$ws.Range("I5").Value = 0.1758510422341793
$ws.Range("J5").Value = 0.1758510422341793
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4770616666666667
$ws.Range("N5").Value = 1.431185
$ws.Range("O5").Value = 0.1014571503512615
$ws.Range("P5").Value = 0.1014571503512615
$ws.Range("Q5").Value = 0.3346050102188889
$ws.Range("R5").Value = 3.011445091970001
$ws.Range("S5").Value = 0.01784134563137917
$ws.Range("T5").Value = 0.01784134563137917

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gphb5"
$ws.Range("C6").Value = "Tshr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.364716666666667
$ws.Range("H6").Value = 7.09415
$ws.Range("I6").Value = 0.5928790992640315
$ws.Range("J6").Value = 0.5928790992640315
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7209850000000001
$ws.Range("N6").Value = 2.162955
$ws.Range("O6").Value = 0.1533325535399077
$ws.Range("P6").Value = 0.1533325535399077
$ws.Range("Q6").Value = 1.704925245916667
$ws.Range("R6").Value = 15.34432721325
$ws.Range("S6").Value = 0.09090766623059439
$ws.Range("T6").Value = 0.09090766623059436

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gphb5"
$ws.Range("C7").Value = "Tshr"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.364716666666667
$ws.Range("H7").Value = 7.09415
$ws.Range("I7").Value = 0.5928790992640315
$ws.Range("J7").Value = 0.5928790992640315
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.387303666666666
$ws.Range("N7").Value = 10.161911
$ws.Range("O7").Value = 0.7203810354238886
$ws.Range("P7").Value = 0.7203810354238885
$ws.Range("Q7").Value = 8.010013435627778
$ws.Range("R7").Value = 72.09012092064999
$ws.Range("S7").Value = 0.4270988594090054
$ws.Range("T7").Value = 0.4270988594090054

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Gphb5"
$ws.Range("C8").Value = "Tshr"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.364716666666667
$ws.Range("H8").Value = 7.09415
$ws.Range("I8").Value = 0.5928790992640315
$ws.Range("J8").Value = 0.5928790992640315
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1167496666666667
$ws.Range("N8").Value = 0.350249
$ws.Range("O8").Value = 0.02482926068494219
$ws.Range("P8").Value = 0.02482926068494219
$ws.Range("Q8").Value = 0.2760798825944445
$ws.Range("R8").Value = 2.48471894335
$ws.Range("S8").Value = 0.01472074971028036
$ws.Range("T8").Value = 0.01472074971028035

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Gphb5"
$ws.Range("C9").Value = "Tshr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.364716666666667
$ws.Range("H9").Value = 7.09415
$ws.Range("I9").Value = 0.5928790992640315
$ws.Range("J9").Value = 0.5928790992640315
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4770616666666667
$ws.Range("N9").Value = 1.431185
$ws.Range("O9").Value = 0.1014571503512615
$ws.Range("P9").Value = 0.1014571503512615
$ws.Range("Q9").Value = 1.128115674194445
$ws.Range("R9").Value = 10.15304106775
$ws.Range("S9").Value = 0.06015182391415134
$ws.Range("T9").Value = 0.06015182391415134

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Gphb5"
$ws.Range("C10").Value = "Tshr"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.922427
$ws.Range("H10").Value = 2.767281
$ws.Range("I10").Value = 0.2312698585017892
$ws.Range("J10").Value = 0.2312698585017893
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7209850000000001
$ws.Range("N10").Value = 2.162955
$ws.Range("O10").Value = 0.1533325535399077
$ws.Range("P10").Value = 0.1533325535399077
$ws.Range("Q10").Value = 0.6650560305950001
$ws.Range("R10").Value = 5.985504275355001
$ws.Range("S10").Value = 0.03546119796089248
$ws.Range("T10").Value = 0.03546119796089248

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Gphb5"
$ws.Range("C11").Value = "Tshr"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.922427
$ws.Range("H11").Value = 2.767281
$ws.Range("I11").Value = 0.2312698585017892
$ws.Range("J11").Value = 0.2312698585017893
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.387303666666666
$ws.Range("N11").Value = 10.161911
$ws.Range("O11").Value = 0.7203810354238886
$ws.Range("P11").Value = 0.7203810354238885
$ws.Range("Q11").Value = 3.124540359332333
$ws.Range("R11").Value = 28.120863233991
$ws.Range("S11").Value = 0.1666024201298552
$ws.Range("T11").Value = 0.1666024201298552

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Gphb5"
$ws.Range("C12").Value = "Tshr"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.922427
$ws.Range("H12").Value = 2.767281
$ws.Range("I12").Value = 0.2312698585017892
$ws.Range("J12").Value = 0.2312698585017893
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1167496666666667
$ws.Range("N12").Value = 0.350249
$ws.Range("O12").Value = 0.02482926068494219
$ws.Range("P12").Value = 0.02482926068494219
$ws.Range("Q12").Value = 0.1076930447743333
$ws.Range("R12").Value = 0.9692374029690001
$ws.Range("S12").Value = 0.00574225960531062
$ws.Range("T12").Value = 0.00574225960531062

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Gphb5"
$ws.Range("C13").Value = "Tshr"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.922427
$ws.Range("H13").Value = 2.767281
$ws.Range("I13").Value = 0.2312698585017892
$ws.Range("J13").Value = 0.2312698585017893
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4770616666666667
$ws.Range("N13").Value = 1.431185
$ws.Range("O13").Value = 0.1014571503512615
$ws.Range("P13").Value = 0.1014571503512615
$ws.Range("Q13").Value = 0.4400545619983334
$ws.Range("R13").Value = 3.960491057985001
$ws.Range("S13").Value = 0.02346398080573101
$ws.Range("T13").Value = 0.02346398080573101
